$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H3").Value = 4.7
$ws.Range("I3").Value = 6.6
$ws.Range("K3").Value = 3.9
$ws.Range("S3").Value = 3.7
$ws.Range("V3").Value = 1.18
$ws.Range("X3").Value = 13.5
$ws.Range("Z3").Value = 48
$ws.Range("F4").Value = 1.77
$ws.Range("U4").Value = 1.74
$ws.Range("AI4").Value = 130
$ws.Range("F5").Value = 7.6
$ws.Range("G5").Value = 7.8
$ws.Range("I5").Value = 1.53
$ws.Range("T5").Value = 1.96
$ws.Range("V5").Value = 2.9
$ws.Range("W5").Value = 1.14
$ws.Range("AL5").Value = 100
$ws.Range("H6").Value = 13
$ws.Range("K6").Value = 7.8
$ws.Range("O6").Value = 1.12
$ws.Range("R6").Value = 1.98
$ws.Range("S6").Value = 1.98
$ws.Range("Z6").Value = 150
$ws.Range("F7").Value = 3.45
$ws.Range("H7").Value = 2.22
$ws.Range("I7").Value = 2.24
$ws.Range("O7").Value = 1.25
$ws.Range("P7").Value = 2.26
$ws.Range("Q7").Value = 1.75
$ws.Range("R7").Value = 1.51
$ws.Range("S7").Value = 2.88
$ws.Range("V7").Value = 1.8
$ws.Range("AE7").Value = 19
$ws.Range("AG7").Value = 14
$ws.Range("AL7").Value = 42
$ws.Range("AO7").Value = 13.5
$ws.Range("G8").Value = 1.46
$ws.Range("H8").Value = 8.199999999999999
$ws.Range("K8").Value = 5.2
$ws.Range("R8").Value = 1.45
$ws.Range("W8").Value = 3.15
$ws.Range("AJ8").Value = 12
$ws.Range("AN8").Value = 7
$ws.Range("Z9").Value = 17.5
$ws.Range("AC9").Value = 9.199999999999999
$ws.Range("AE9").Value = 20
$ws.Range("AH9").Value = 14.5
$ws.Range("AN9").Value = 20
$ws.Range("AO10").Value = 18
$ws.Range("G11").Value = 2.22
$ws.Range("K11").Value = 3.85
$ws.Range("S11").Value = 2.78
$ws.Range("X11").Value = 18
$ws.Range("AC11").Value = 8.800000000000001
$ws.Range("AG11").Value = 10.5
$ws.Range("AK11").Value = 20
$ws.Range("AN11").Value = 12.5
$ws.Range("H12").Value = 18.5
$ws.Range("I12").Value = 19
$ws.Range("P12").Value = 4.3
$ws.Range("T12").Value = 1.83
$ws.Range("AF12").Value = 12
$ws.Range("AN12").Value = 2.46
$ws.Range("P13").Value = 2.32
$ws.Range("R13").Value = 1.53
